$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'48.130.88"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +1.71%  "
$ws.Range("D3").Value = "'2.507.16"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +0.60%  "
$ws.Range("E4").Value = "  +0.05%  "
$ws.Range("D5").Value = "'109.06"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.43%  "
$ws.Range("D6").Value = "'320.52"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.13%  "
$ws.Range("D7").Value = "'0.529"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +1.28%  "
$ws.Range("E9").Value = "  +1.71%  "
$ws.Range("D10").Value = "'39.82"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +1.69%  "
$ws.Range("D11").Value = "'20.09"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +9.29%  "
$ws.Range("D12").Value = "'0.0817"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +0.88%  "
$ws.Range("D13").Value = "'0.124"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +0.68%  "
$ws.Range("D14").Value = "'7.19"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +0.81%  "
$ws.Range("D15").Value = "'2.903.93"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +0.78%  "
$ws.Range("D16").Value = "'2.510.65"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +0.56%  "
$ws.Range("E17").Value = "  -0.05%  "
$ws.Range("D18").Value = "'47.985.03"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +1.60%  "
$ws.Range("D19").Value = "'13.17"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.75%  "
$ws.Range("D20").Value = "'6.62"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +0.06%  "
$ws.Range("D21").Value = "'0.0₃0941"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.67%  "
$ws.Range("D22").Value = "'2.72"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +2.51%  "
$ws.Range("D23").Value = "'72.22"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +2.60%  "
$ws.Range("D24").Value = "'274.75"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +12.08%  "
$ws.Range("D25").Value = "'2.57"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.35%  "
$ws.Range("E26").Value = "  -0.02%  "
$ws.Range("D27").Value = "'25.89"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +0.52%  "
$ws.Range("D28").Value = "'2.40"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +7.43%  "
$ws.Range("D29").Value = "'10.07"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +0.59%  "
$ws.Range("E30").Value = "  +2.38%  "
$ws.Range("D31").Value = "'35.54"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +2.26%  "
$ws.Range("D32").Value = "'49.45"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -0.52%  "
$ws.Range("D33").Value = "'19.29"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -6.42%  "
$ws.Range("E34").Value = "  -0.18%  "
$ws.Range("E35").Value = "  -0.01%  "
$ws.Range("E36").Value = "  +0.03%  "
$ws.Range("E37").Value = "  -0.46%  "
$ws.Range("D38").Value = "'4.62"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -2.30%  "
$ws.Range("E39").Value = "  +0.83%  "
$ws.Range("E40").Value = "  +0.85%  "
$ws.Range("D41").Value = "'122.62"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +5.74%  "
$ws.Range("E42").Value = "  -0.62%  "
$ws.Range("D43").Value = "'21.84"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -6.28%  "
$ws.Range("D44").Value = "'0.0304"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +2.60%  "
$ws.Range("D45").Value = "'2.024.48"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +1.40%  "
$ws.Range("E46").Value = "  +3.11%  "
$ws.Range("D47").Value = "'1.85"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +3.88%  "
$ws.Range("E48").Value = "  -0.69%  "
$ws.Range("D49").Value = "'9.03"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -1.15%  "
$ws.Range("D50").Value = "'5.18"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +1.68%  "
$ws.Range("D51").Value = "'79.35"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +2.60%  "
